# Daily attendance processing - 2026-01-29 22:42:34
# Reorders the comma-separated "Recorded By" list (column G) for a fixed
# set of exact values, moving the trailing item to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cur = $cell.Text
    if ($map.ContainsKey($cur)) {
        $cell.Value = $map[$cur]
    }
}
